$wb = $excel.ActiveWorkbook

# --- "Veicolo" sheet: was empty, now gets a small "id" list ---
$wsVeicolo = $wb.Worksheets.Item("Veicolo")
$wsVeicolo.Range("A1").Value = "id"
$wsVeicolo.Range("A2").Value = 252681
$wsVeicolo.Range("A3").Value = 253497
$wsVeicolo.Columns.Item(1).ColumnWidth = 29.2

# --- "Release Date (RD)" sheet: previously had an id/release_date/tassativita
#     table with 3 rows; now cleared out entirely (no RD issues found) ---
$wsRD = $wb.Worksheets.Item("Release Date (RD)")
$wsRD.Range("A1:C4").EntireRow.Delete()

# --- "RD Tassative" sheet: previously had 15 data rows; now reduced to a
#     single row with a new id/date pair ---
$wsRDT = $wb.Worksheets.Item("RD Tassative")
$wsRDT.Range("A3:B16").EntireRow.Delete()
$wsRDT.Range("A2").Value = 253295
$wsRDT.Range("B2").Value = 45911.58333333334
